$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the new columns are formatted as text so values like "306.00"
# or "75.90" keep their literal representation instead of being coerced
# to numbers (matching the original sheet, where every cell is inlineStr).
$ws.Columns("Z:AH").NumberFormat = "@"

# Row 1 - time headers
$ws.Range("Z1").Value = "13:18"
$ws.Range("AA1").Value = "13:19"
$ws.Range("AB1").Value = "13:20"
$ws.Range("AC1").Value = "13:24"
$ws.Range("AD1").Value = "13:25"
$ws.Range("AE1").Value = "13:26"
$ws.Range("AF1").Value = "13:27"
$ws.Range("AG1").Value = "13:29"
$ws.Range("AH1").Value = "13:31"

# Row 2 - HUBC
$ws.Range("Z2").Value = "75.99"
$ws.Range("AA2").Value = "75.99"
$ws.Range("AB2").Value = "75.86"
$ws.Range("AC2").Value = "75.86"
$ws.Range("AD2").Value = "75.98"
$ws.Range("AE2").Value = "75.98"
$ws.Range("AF2").Value = "75.98"
$ws.Range("AG2").Value = "75.98"
$ws.Range("AH2").Value = "75.90"

# Row 3 - GATI
$ws.Range("Z3").Value = "306.00"
$ws.Range("AA3").Value = "306.00"
$ws.Range("AB3").Value = "306.00"
$ws.Range("AC3").Value = "306.00"
$ws.Range("AD3").Value = "306.00"
$ws.Range("AE3").Value = "306.00"
$ws.Range("AF3").Value = "306.00"
$ws.Range("AG3").Value = "306.00"
$ws.Range("AH3").Value = "306.00"

# Row 4 - OGDC
$ws.Range("Z4").Value = "147.35"
$ws.Range("AA4").Value = "147.35"
$ws.Range("AB4").Value = "147.26"
$ws.Range("AC4").Value = "147.26"
$ws.Range("AD4").Value = "147.25"
$ws.Range("AE4").Value = "147.25"
$ws.Range("AF4").Value = "147.25"
$ws.Range("AG4").Value = "147.25"
$ws.Range("AH4").Value = "147.02"
